$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price values look numeric must be forced to Text format
# first, otherwise Excel will auto-convert the string into a real number and
# the literal text (e.g. "1.009") would be lost/rounded.
$textPriceCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D14", "D15", "D16", "D17", "D18", "D20", "D22", "D25", "D26", "D27", "D28", "D30", "D31", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50")
foreach ($cellAddr in $textPriceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Update Price (column D) and Volume(1h) (column E) values row by row
$ws.Range("D2").Value = "27.651.07"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.845.04"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  -1.85%  "
$ws.Range("D5").Value = "317.68"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("D7").Value = "0.4297"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Value = "0.07318"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").Value = "0.8729"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "1.844.95"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "5.420"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "0.07114"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "88.67"
$ws.Range("E16").Value = "  +4.22%  "
$ws.Range("D17").Value = "1.012"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "0.000008994"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D20").Value = "15.40"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "27.665.52"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "5.200"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("D24").Value = "2.072.06"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").Value = "1.968"
$ws.Range("E25").Value = "  -3.21%  "
$ws.Range("D26").Value = "154.93"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "18.58"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "2.159"
$ws.Range("E28").Value = "  +8.07%  "
$ws.Range("D30").Value = "118.04"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "0.08907"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "0.7739"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").Value = "4.527"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "2.886"
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("D37").Value = "1.128"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D38").Value = "0.01976"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").Value = "0.05320"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").Value = "2.884"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "7.130"
$ws.Range("D42").Value = "0.1689"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("D43").Value = "0.5120"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").Value = "8.767"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "10.69"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").Value = "107.35"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").Value = "0.06446"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D49").Value = "1.010"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").Value = "1.687"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("E51").Value = "  -2.77%  "
